# [LPF-879]: CCMS Third party report
# Remove the sheets/pivot tables that are no longer part of this report,
# keeping only the "Summary" sheet and the "Transparency Rec" sheet
# (with its "Transparency rec pivot" PivotTable / shared pivot cache).

$wb = $excel.ActiveWorkbook

# Suppress the "delete sheet" confirmation dialog.
$excel.DisplayAlerts = $false

$sheetsToRemove = @(
    "By Source and Expenditure type",
    "Provider Contigency",
    "MAIN"
)

foreach ($sheetName in $sheetsToRemove) {
    $wb.Worksheets($sheetName).Delete()
}

$excel.DisplayAlerts = $true
